{"js": "// Update the date heading paragraph (first paragraph in the document body)\nconst headingParagraphs = context.document.body.paragraphs;\nheadingParagraphs.load(\"items\");\nawait context.sync();\nheadingParagraphs.items[0].insertText(\"2025-12-03 Wednesday\", Word.InsertLocation.replace);\n\n// Replacement values for the 20x5 practice table, in row-major order\nconst values = [\n  \"40-34=\",\n  \"8+76=\",\n  \"44+9=\",\n  \"18+38=\",\n  \"20-14=\",\n  \"56-19=\",\n  \"27+55=\",\n  \"16+15=\",\n  \"64+9=\",\n  \"58+38=\",\n  \"59+18=\",\n  \"73-69=\",\n  \"81-18=\",\n  \"80-64=\",\n  \"16+8=\",\n  \"57+9=\",\n  \"52-19=\",\n  \"94-79=\",\n  \"60-49=\",\n  \"14+19=\",\n  \"43-6=\",\n  \"28+4=\",\n  \"54-36=\",\n  \"64-48=\",\n  \"90-81=\",\n  \"93-65=\",\n  \"48+5=\",\n  \"8+73=\",\n  \"55-46=\",\n  \"39+28=\",\n  \"91-53=\",\n  \"37+45=\",\n  \"57-8=\",\n  \"9+82=\",\n  \"82-25=\",\n  \"62-13=\",\n  \"80-67=\",\n  \"16+27=\",\n  \"82-55=\",\n  \"66-8=\",\n  \"48+44=\",\n  \"67-59=\",\n  \"80-63=\",\n  \"95-57=\",\n  \"95-48=\",\n  \"60-51=\",\n  \"16+56=\",\n  \"70-44=\",\n  \"46+18=\",\n  \"7+58=\",\n  \"3+19=\",\n  \"4+89=\",\n  \"27+38=\",\n  \"91-6=\",\n  \"75-16=\",\n  \"25+36=\",\n  \"23+28=\",\n  \"67+29=\",\n  \"7+26=\",\n  \"30-16=\",\n  \"19+4=\",\n  \"51-24=\",\n  \"52+9=\",\n  \"70-47=\",\n  \"92-33=\",\n  \"37+45=\",\n  \"10-5=\",\n  \"37+36=\",\n  \"82-79=\",\n  \"30-6=\",\n  \"17+29=\",\n  \"35+26=\",\n  \"29+63=\",\n  \"82-9=\",\n  \"38+58=\",\n  \"91-19=\",\n  \"43-19=\",\n  \"2+39=\",\n  \"24+47=\",\n  \"39+39=\",\n  \"29+32=\",\n  \"91-2=\",\n  \"29+34=\",\n  \"62-33=\",\n  \"8+35=\",\n  \"81-55=\",\n  \"74-68=\",\n  \"6+9=\",\n  \"93-66=\",\n  \"8+57=\",\n  \"71-37=\",\n  \"7+27=\",\n  \"81-37=\",\n  \"13+78=\",\n  \"6+8=\",\n  \"68+6=\",\n  \"84+8=\",\n  \"86-19=\",\n  \"34+8=\",\n  \"87-68=\"\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nlet idx = 0;\nfor (let r = 0; r < 20; r++) {\n  for (let c = 0; c < 5; c++) {\n    const cell = table.getCell(r, c);\n    cell.value = values[idx];\n    idx++;\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the date heading paragraph\n$p = $d.Paragraphs.Item(1)\n$p.Range.Text = \"2025-12-03 Wednesday\"\n\n# Replacement values for the 20x5 practice table, in row-major order\n$values = @(\n  \"40-34=\",\n  \"8+76=\",\n  \"44+9=\",\n  \"18+38=\",\n  \"20-14=\",\n  \"56-19=\",\n  \"27+55=\",\n  \"16+15=\",\n  \"64+9=\",\n  \"58+38=\",\n  \"59+18=\",\n  \"73-69=\",\n  \"81-18=\",\n  \"80-64=\",\n  \"16+8=\",\n  \"57+9=\",\n  \"52-19=\",\n  \"94-79=\",\n  \"60-49=\",\n  \"14+19=\",\n  \"43-6=\",\n  \"28+4=\",\n  \"54-36=\",\n  \"64-48=\",\n  \"90-81=\",\n  \"93-65=\",\n  \"48+5=\",\n  \"8+73=\",\n  \"55-46=\",\n  \"39+28=\",\n  \"91-53=\",\n  \"37+45=\",\n  \"57-8=\",\n  \"9+82=\",\n  \"82-25=\",\n  \"62-13=\",\n  \"80-67=\",\n  \"16+27=\",\n  \"82-55=\",\n  \"66-8=\",\n  \"48+44=\",\n  \"67-59=\",\n  \"80-63=\",\n  \"95-57=\",\n  \"95-48=\",\n  \"60-51=\",\n  \"16+56=\",\n  \"70-44=\",\n  \"46+18=\",\n  \"7+58=\",\n  \"3+19=\",\n  \"4+89=\",\n  \"27+38=\",\n  \"91-6=\",\n  \"75-16=\",\n  \"25+36=\",\n  \"23+28=\",\n  \"67+29=\",\n  \"7+26=\",\n  \"30-16=\",\n  \"19+4=\",\n  \"51-24=\",\n  \"52+9=\",\n  \"70-47=\",\n  \"92-33=\",\n  \"37+45=\",\n  \"10-5=\",\n  \"37+36=\",\n  \"82-79=\",\n  \"30-6=\",\n  \"17+29=\",\n  \"35+26=\",\n  \"29+63=\",\n  \"82-9=\",\n  \"38+58=\",\n  \"91-19=\",\n  \"43-19=\",\n  \"2+39=\",\n  \"24+47=\",\n  \"39+39=\",\n  \"29+32=\",\n  \"91-2=\",\n  \"29+34=\",\n  \"62-33=\",\n  \"8+35=\",\n  \"81-55=\",\n  \"74-68=\",\n  \"6+9=\",\n  \"93-66=\",\n  \"8+57=\",\n  \"71-37=\",\n  \"7+27=\",\n  \"81-37=\",\n  \"13+78=\",\n  \"6+8=\",\n  \"68+6=\",\n  \"84+8=\",\n  \"86-19=\",\n  \"34+8=\",\n  \"87-68=\"\n)\n\n$t = $d.Tables.Item(1)\n$idx = 0\nfor ($r = 1; $r -le 20; $r++) {\n    for ($c = 1; $c -le 5; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $values[$idx]\n        $idx = $idx + 1\n    }\n}\n\nWrite-Output \"Done: updated $idx cells\"\n"}
